# Update cryptos list with latest scraped prices / volume changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.028.96"
$ws.Range("E2").Value = "  -0.37%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.643.83"
$ws.Range("E3").Value = "  +0.10%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.75%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.42"
$ws.Range("E5").Value = "  -0.10%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.79%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.54%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.15%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +0.60%  "

# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.63"
$ws.Range("E10").Value = "  -0.44%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.66%  "

# Row 12 - Polkadot
$ws.Range("E12").Value = "  +0.19%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.869.78"
$ws.Range("E13").Value = "  +0.06%  "

# Row 14 - WrappedEther
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.631.95"
$ws.Range("E14").Value = "  -0.16%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.27%  "

# Row 17 - Litecoin
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.05"
$ws.Range("E17").Value = "  -0.32%  "

# Row 18 - WrappedBTC
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.942.35"
$ws.Range("E18").Value = "  -0.67%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.59%  "

# Row 20 - BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.18"
$ws.Range("E20").Value = "  -1.02%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -1.51%  "

# Row 22 - Avalanche
$ws.Range("E22").Value = "  -0.77%  "

# Row 23 - Chainlink
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.27"
$ws.Range("E23").Value = "  -0.07%  "

# Row 24 - now Stellar (was Toncoin)
$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.132"
$ws.Range("E24").Value = "  +6.38%  "

# Row 25 - now Toncoin (was Stellar)
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.81"
$ws.Range("E25").Value = "  +1.20%  "

# Row 26 - Monero
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.84"
$ws.Range("E26").Value = "  +1.56%  "

# Row 27 - BinanceUSD
$ws.Range("E27").Value = "  +0.63%  "

# Row 28 - Cosmos
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.94"
$ws.Range("E28").Value = "  +0.31%  "

# Row 29 - EthereumClassic
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.58"
$ws.Range("E29").Value = "  -0.04%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.21%  "

# Row 31 - Hedera
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0501"
$ws.Range("E31").Value = "  -0.14%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.30"
$ws.Range("E32").Value = "  -1.68%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  +0.40%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  -3.21%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  +2.37%  "

# Row 36 - ARBITRUM
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.906"
$ws.Range("E36").Value = "  -0.62%  "

# Row 37 - Maker
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.135.45"
$ws.Range("E37").Value = "  -0.01%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  -1.43%  "

# Row 39 - MXToken
$ws.Range("E39").Value = "  -1.42%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +0.34%  "

# Row 41 - FraxShare
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.52"
$ws.Range("E41").Value = "  +0.71%  "

# Row 42 - Quant
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.70"
$ws.Range("E42").Value = "  -0.43%  "

# Row 43 - TrustWalletToken
$ws.Range("E43").Value = "  +0.20%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.779.23"
$ws.Range("E44").Value = "  +0.06%  "

# Row 45 - BabyDogeCoin
$ws.Range("E45").Value = "  +3.30%  "

# Row 46 - Aave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.87"
$ws.Range("E46").Value = "  +0.11%  "

# Row 47 - Cronos
$ws.Range("E47").Value = "  +2.74%  "

# Row 48 - RenderToken
$ws.Range("E48").Value = "  -0.84%  "

# Row 49 - EnergySwap
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.74"
$ws.Range("E49").Value = "  +1.12%  "

# Row 50 - Mantle
$ws.Range("E50").Value = "  -0.20%  "

# Row 51 - Algorand
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0963"
$ws.Range("E51").Value = "  -0.14%  "
